$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 14: new date header block (mirrors row 1's "ago_force_voltage_slope" /
# "ant_force_voltage_slope" header row, but for the 2024/08/11 data set)
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 45515
$ws.Range("D14").Value = "ago_force_voltage_slope"
$ws.Range("E14").Value = "ant_force_voltage_slope"

# Rows 15-19: raw samples
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 98.836269999999999

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = -83.661230000000003
$ws.Range("E16").Value = 98.916240000000002

$ws.Range("C17").Value = 3
$ws.Range("D17").Value = -67.876499999999993
$ws.Range("E17").Value = 76.397530000000003

$ws.Range("C18").Value = 4
$ws.Range("D18").Value = -87.812299999999993
$ws.Range("E18").Value = 73.182720000000003

$ws.Range("C19").Value = 5
$ws.Range("D19").Value = -126.319

# Row 20: averages
$ws.Range("C20").Value = "ave"
$ws.Range("D20").Formula = "=AVERAGE(D15:D19)"
$ws.Range("E20").Formula = "=AVERAGE(E15:E18)"

# ---------------------------------------------------------------------------
# Rows 23-30: recompute the a_n_agonist / a_n_antagonist slopes against the
# new row-20 averages (same layout as rows 4-11, new anchor row)
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "a_3_agonist"
$ws.Range("B23").Value = -0.20082
$ws.Range("D23").Formula = "=B23*`$D`$20"

$ws.Range("A24").Value = "a_2_agonist"
$ws.Range("B24").Value = 7.001792
$ws.Range("A25").Value = "a_1_agonist"
$ws.Range("B25").Value = 0.25617299999999998
$ws.Range("A26").Value = "a_0_agonist"
$ws.Range("B26").Value = 0.91102000000000005
$ws.Range("D24:D26").Formula = "=B24*`$D`$20"

$ws.Range("A27").Value = "a_3_antagonist"
$ws.Range("B27").Value = 0.245392
$ws.Range("E27").Formula = "=B27*`$E`$20"

$ws.Range("A28").Value = "a_2_antagonist"
$ws.Range("B28").Value = 0.34852100000000003
$ws.Range("A29").Value = "a_1_antagonist"
$ws.Range("B29").Value = -0.036479999999999999
$ws.Range("A30").Value = "a_0_antagonist"
$ws.Range("B30").Value = 6.3177130000000004
$ws.Range("E28:E30").Formula = "=B28*`$E`$20"

# ---------------------------------------------------------------------------
# Row 32: section header for the direct force-gauge re-measurement
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = "2024/08/12（直接フォースゲージつけた）"

# ---------------------------------------------------------------------------
# Rows 33-40: agonist/antagonist coefficients (A/B) plus the two new
# force/voltage ratio tables (D:F for ago, H:J for anta) and the resulting
# L/M slope columns.
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = "a_3_agonist"
$ws.Range("B33").Value = -0.20082
$ws.Range("D33").Value = "ago_force"
$ws.Range("E33").Value = "ago_voltage"
$ws.Range("H33").Value = "anta_force"
$ws.Range("I33").Value = "anta_voltage"
$ws.Range("L33").Formula = "=B33*`$F`$39"

$ws.Range("A34").Value = "a_2_agonist"
$ws.Range("B34").Value = 7.001792
$ws.Range("C34").Value = 1
$ws.Range("H34").Value = 11.9
$ws.Range("I34").Value = 39
$ws.Range("J34").Formula = "=I34/H34"
$ws.Range("L34").Formula = "=B34*`$F`$39"

$ws.Range("A35").Value = "a_1_agonist"
$ws.Range("B35").Value = 0.25617299999999998
$ws.Range("C35").Value = 2
$ws.Range("D35").Value = 22.9
$ws.Range("E35").Value = -48
$ws.Range("F35").Formula = "=E35/D35"
$ws.Range("L35").Formula = "=B35*`$F`$39"

$ws.Range("A36").Value = "a_0_agonist"
$ws.Range("B36").Value = 0.91102000000000005
$ws.Range("C36").Value = 3
$ws.Range("D36").Value = 20.6
$ws.Range("E36").Value = -27
$ws.Range("H36").Value = 10.9
$ws.Range("I36").Value = 34
$ws.Range("J36").Formula = "=I36/H36"
$ws.Range("L36").Formula = "=B36*`$F`$39"

$ws.Range("A37").Value = "a_3_antagonist"
$ws.Range("B37").Value = 0.245392
$ws.Range("C37").Value = 4
$ws.Range("D37").Value = 18.5
$ws.Range("E37").Value = -24
$ws.Range("H37").Value = 12.4
$ws.Range("I37").Value = 45
$ws.Range("J37").Formula = "=I37/H37"
$ws.Range("M37").Formula = "=B37*`$J`$39"

$ws.Range("A38").Value = "a_2_antagonist"
$ws.Range("B38").Value = 0.34852100000000003
$ws.Range("C38").Value = 5
$ws.Range("D38").Value = 23.8
$ws.Range("E38").Value = -63
$ws.Range("H38").Value = 12
$ws.Range("I38").Value = 50
$ws.Range("J38").Formula = "=I38/H38"

$ws.Range("A39").Value = "a_1_antagonist"
$ws.Range("B39").Value = -0.036479999999999999

$ws.Range("A40").Value = "a_0_antagonist"
$ws.Range("B40").Value = 6.3177130000000004

# F36:F38 and M38:M40 are contiguous blocks -> fill as a range so Excel
# records them as a single shared formula, matching the source fill-handle
# drag. The F/L/J/M formulas above (with gaps) are set cell-by-cell since
# their "ref" spans in the source aren't fully populated.
$ws.Range("F36:F38").Formula = "=E36/D36"
$ws.Range("M38:M40").Formula = "=B38*`$J`$39"

$ws.Range("F39").Formula = "=AVERAGE(F35:F38)"
$ws.Range("J39").Formula = "=AVERAGE(J34:J38)"

# ---------------------------------------------------------------------------
# View / print setup
# ---------------------------------------------------------------------------
$ws.Range("E19").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("M37:M40").Select()

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
